$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "31.198.65"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +2.16%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.939.68"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.42%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "241.93"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.95%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06781"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "20.16"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +7.30%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "104.13"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07841"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.24%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.947.78"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.38%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.296"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6983"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "295.90"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +11.28%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "31.189.52"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.98%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.208.97"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("E19").Value = "  +2.15%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007612"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.14%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.570"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("E22").Value = "  -0.28%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.40%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.420"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.560"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "169.36"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.41%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "19.81"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +3.82%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.101"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.394"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.37%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.1008"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.68%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.627"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +1.53%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.04839"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.7391"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.135"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.724"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.20%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01960"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +0.60%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.862"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +8.05%  "
$ws.Range("E40").Value = "  -0.39%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "76.56"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +1.26%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.8727"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.13%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.4376"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.70%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "105.89"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").Value = "  -0.38%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.028.14"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.99%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.589"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.312"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +3.66%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.1209"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "35.24"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
